# Natmi following Dr Hou advice
# Re-generated NATMI LR-pair output: sending/target clusters now include the
# "ECs" cluster as a sending cluster too (full 3x3 grid of FAPs/sCs/ECs),
# so shared-string order changes (ECs moves before FAPs) and every data row
# is recomputed; 3 new rows are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Sema3b", "Nrp2", "ECs", 2, 0.6666666666666666, 2.025829666666667, 6.077489, 0.1667292296771927, 0.1667292296771927, 3, 1, 34.52052333333333, 103.56157, 0.7684334662422598, 0.7684334662422598, 69.93270027752554, 629.3943024977299, 0.128120319884747, 0.128120319884747),
    @("ECs", "Sema3b", "Nrp2", "FAPs", 2, 0.6666666666666666, 2.025829666666667, 6.077489, 0.1667292296771927, 0.1667292296771927, 3, 1, 5.347618333333333, 16.042855, 0.1190390091234806, 0.1190390091234805, 10.83336386567722, 97.50027479109499, 0.01984728229269423, 0.01984728229269422),
    @("ECs", "Sema3b", "Nrp2", "sCs", 2, 0.6666666666666666, 2.025829666666667, 6.077489, 0.1667292296771927, 0.1667292296771927, 3, 1, 5.055101333333334, 15.165304, 0.1125275246342597, 0.1125275246342597, 10.24077424907289, 92.16696824165601, 0.01876162749975144, 0.01876162749975144),
    @("FAPs", "Sema3b", "Nrp2", "ECs", 3, 1, 5.787456, 17.362368, 0.4763174794741613, 0.4763174794741613, 3, 1, 34.52052333333333, 103.56157, 0.7684334662422598, 0.7684334662422598, 199.78600988864, 1798.07408899776, 0.3660182917841062, 0.3660182917841062),
    @("FAPs", "Sema3b", "Nrp2", "FAPs", 3, 1, 5.787456, 17.362368, 0.4763174794741613, 0.4763174794741613, 3, 1, 5.347618333333333, 16.042855, 0.1190390091234806, 0.1190390091234805, 30.94910580896, 278.54195228064, 0.05670036078479795, 0.05670036078479794),
    @("FAPs", "Sema3b", "Nrp2", "sCs", 3, 1, 5.787456, 17.362368, 0.4763174794741613, 0.4763174794741613, 3, 1, 5.055101333333334, 15.165304, 0.1125275246342597, 0.1125275246342597, 29.256176542208, 263.305588879872, 0.05359882690525717, 0.05359882690525716),
    @("sCs", "Sema3b", "Nrp2", "ECs", 3, 1, 4.337131333333333, 13.011394, 0.356953290848646, 0.356953290848646, 3, 1, 34.52052333333333, 103.56157, 0.7684334662422598, 0.7684334662422598, 149.7200433920644, 1347.48039052858, 0.2742948545734066, 0.2742948545734066),
    @("sCs", "Sema3b", "Nrp2", "FAPs", 3, 1, 4.337131333333333, 13.011394, 0.356953290848646, 0.356953290848646, 3, 1, 5.347618333333333, 16.042855, 0.1190390091234806, 0.1190390091234805, 23.19332303220778, 208.73990728987, 0.04249136604598838, 0.04249136604598838),
    @("sCs", "Sema3b", "Nrp2", "sCs", 3, 1, 4.337131333333333, 13.011394, 0.356953290848646, 0.356953290848646, 3, 1, 5.055101333333334, 15.165304, 0.1125275246342597, 0.1125275246342597, 21.92463838597511, 197.321745473776, 0.04016707022925108, 0.04016707022925108)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $row[$c]
    }
}
